# Weekly update: insert a new price record as the new row 133 (most recent
# week for this Terminal Hortofrutícola Agro Chillán - Zapallo italiano
# series) and push the previously-existing rows 133:139 down to 134:140.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 133:139 down to 134:140, leaving a blank row 133 to populate.
$ws.Rows("133:133").Insert()

# Populate the new row 133 with this week's data.
$ws.Cells.Item(133, 1).Value  = 7
$ws.Cells.Item(133, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(133, 3).Value  = "Ñuble"
$ws.Cells.Item(133, 4).Value  = 44509
$ws.Cells.Item(133, 5).Value  = 16
$ws.Cells.Item(133, 6).Value  = 100112032
$ws.Cells.Item(133, 7).Value  = "Zapallo italiano"
$ws.Cells.Item(133, 8).Value  = "Sin especificar"
$ws.Cells.Item(133, 9).Value  = "Primera"
$ws.Cells.Item(133, 10).Value = 100
$ws.Cells.Item(133, 11).Value = 14000
$ws.Cells.Item(133, 12).Value = 15000
$ws.Cells.Item(133, 13).Value = 14500
$ws.Cells.Item(133, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(133, 15).Value = "Región del Maule"
$ws.Cells.Item(133, 16).Value = 242
$ws.Cells.Item(133, 17).Value = 60
$ws.Cells.Item(133, 18).Value = "Hortaliza"
